$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.235.20"
$ws.Range("E2").Value = "  -3.19%  "
$ws.Range("D3").Value = "1.926.01"
$ws.Range("E3").Value = "  -2.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7152"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3244"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.35"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06820"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8002"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07932"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "1.926.16"
$ws.Range("E13").Value = "  -2.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.393"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "260.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.35%  "
$ws.Range("D18").Value = "30.240.26"
$ws.Range("E18").Value = "  -3.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007939"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.811"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "2.178.81"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.852"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.626"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.80%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1329"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -12.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.280"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.358"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.548"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.417"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.186"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05062"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.191"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7386"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.739"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01934"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.810"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.527"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4446"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.8304"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.684"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.254"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.474"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.22%  "
